$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one "Perejil" price record per row (rows 2..157), each row
# carrying a reporting date (column D) plus price figures (J, K, L, M, P).
# A new weekly record is inserted at row 72; every existing record from the
# old row 72 onward is pushed down by one row, and the last existing record
# (old row 157) lands in a brand-new row 158.

# Step 1: append the new row 158 by duplicating row 157 (boilerplate columns
# A,B,C,E,F,G,H,I,N,O,Q,R plus the still-unshifted D/J/K/L/M/P values, which
# is exactly what the new last row should contain).
$ws.Range("A157:R157").Copy($ws.Range("A158:R158"))

# Step 2: shift the date/price columns down by one row, working from the
# bottom (157) up to (73) so each source row is read before it gets
# overwritten.
for ($r = 157; $r -ge 73; $r--) {
    $prev = $r - 1
    $ws.Range("D$r").Value2 = $ws.Range("D$prev").Value2
    $ws.Range("J$r").Value2 = $ws.Range("J$prev").Value2
    $ws.Range("K$r").Value2 = $ws.Range("K$prev").Value2
    $ws.Range("L$r").Value2 = $ws.Range("L$prev").Value2
    $ws.Range("M$r").Value2 = $ws.Range("M$prev").Value2
    $ws.Range("P$r").Value2 = $ws.Range("P$prev").Value2
}

# Step 3: write the brand-new record into row 72.
$ws.Range("D72").Value2 = 44763
$ws.Range("J72").Value2 = 2000
$ws.Range("K72").Value2 = 2000
$ws.Range("L72").Value2 = 2500
$ws.Range("M72").Value2 = 2250
$ws.Range("P72").Value2 = 1500
